$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44 and 45 swap content (Maker <-> FirstDigitalUSD) plus updated price/volume values
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.070.01"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "

# Price (D) and Volume(1h) (E) updates for all other rows
$ws.Range("D2").Value = "42.732.30"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "2.519.73"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.579"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0810"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.54"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").Value = "2.906.19"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.28%  "
$ws.Range("D16").Value = "2.513.94"
$ws.Range("E16").Value = "  -4.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.863"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("D18").Value = "42.731.98"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.32%  "
$ws.Range("D20").Value = "0.0₃0969"
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("E25").Value = "  -2.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.94%  "
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("E28").Value = "  +12.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.64%  "
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("E35").Value = "  -3.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0786"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.44%  "
$ws.Range("E37").Value = "  -5.32%  "
$ws.Range("E38").Value = "  -1.84%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.52%  "
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0301"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.13%  "
$ws.Range("D49").Value = "2.759.41"
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("E51").Value = "  -0.44%  "
